$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1) - column F updates ("想去人数")
$ws1.Range("F2").Value  = 2835
$ws1.Range("F3").Value  = 1587
$ws1.Range("F6").Value  = 9640
$ws1.Range("F10").Value = 271
$ws1.Range("F13").Value = 701
$ws1.Range("F14").Value = 701
$ws1.Range("F16").Value = 1206
$ws1.Range("F19").Value = 2264
$ws1.Range("F21").Value = 1964
$ws1.Range("F26").Value = 310
$ws1.Range("F34").Value = 523
$ws1.Range("F39").Value = 1514
$ws1.Range("F40").Value = 43
$ws1.Range("F43").Value = 382
$ws1.Range("F44").Value = 784

# Sheet "全部类型" (sheet4) - column F updates ("想去人数")
$ws4.Range("F2").Value  = 2835
$ws4.Range("F3").Value  = 1587
$ws4.Range("F6").Value  = 9640
$ws4.Range("F12").Value = 271
$ws4.Range("F15").Value = 701
$ws4.Range("F16").Value = 701
$ws4.Range("F17").Value = 1206
$ws4.Range("F20").Value = 2264
$ws4.Range("F21").Value = 1964
$ws4.Range("F25").Value = 310
$ws4.Range("F33").Value = 523
$ws4.Range("F42").Value = 1514
$ws4.Range("F43").Value = 43
$ws4.Range("F47").Value = 382
$ws4.Range("F48").Value = 784
